$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert 5 new rows before the current row 110 (old "bk_status" row),
# pushing the existing rows 110:145 down to 115:150.
$ws.Range("A110:A114").EntireRow.Insert()

# --- Populate the 5 new rows with the new svc_sub_* field definitions ---

# Row 110: svc_sub_name
$ws.Range("A110").Value = "svc_sub_name"
$ws.Range("B110").Value = "Categorical"
$ws.Range("C110").Value = "strs"
$ws.Range("D110").Value = "str"

# Row 111: svc_sub_fee_type
$ws.Range("A111").Value = "svc_sub_fee_type"
$ws.Range("B111").Value = "Categorical"
$ws.Range("C111").Value = "strs"
$ws.Range("D111").Value = "str"

# Row 112: svc_sub_fee_gross
$ws.Range("A112").Value = "svc_sub_fee_gross"
$ws.Range("B112").Value = "Numeric"
$ws.Range("C112").Value = "floats"
$ws.Range("D112").Value = "np.float64"

# Row 113: svc_sub_fee2_type
$ws.Range("A113").Value = "svc_sub_fee2_type"
$ws.Range("B113").Value = "Categorical"
$ws.Range("C113").Value = "strs"
$ws.Range("D113").Value = "str"

# Row 114: svc_sub_fee2_gross
$ws.Range("A114").Value = "svc_sub_fee2_gross"
$ws.Range("B114").Value = "Numeric"
$ws.Range("C114").Value = "floats"
$ws.Range("D114").Value = "np.float64"

# Highlight the new field-name cells (same yellow fill used for other
# recently-added rows further down the sheet, cellXf index 17).
$ws.Range("A110:A114").Interior.Color = 65535

# New rows only carry data through column E (which stays empty, matching
# the blank/styled E column used on neighboring rows).
$ws.Range("E110:E114").Value = ""

# --- Refresh the AutoFilter + _FilterDatabase defined name over the grown range ---
$ws.AutoFilterMode = $false
$ws.Range("A1:X150").AutoFilter()

foreach ($n in $wb.Names) {
  if ($n.Name -eq "Sheet1!_FilterDatabase") {
    $n.RefersTo = "=Sheet1!`$A`$1:`$X`$150"
  }
}

# --- View-state tweaks from the edit: scroll/selection on Sheet1 ---
$ws.Activate()
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$ws.Range("B56").Select()
$win.FreezePanes = $true
$ws.Range("E99").Select()
